# Generate Report for Handoff
# Replaces the old localization file UUID/hash with new ones, refreshes the
# handoff/handback timestamps, and clears the stale "Latest Target
# File"/"Latest Handback File" values (and their hyperlinks) on the
# per-language sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "fa571b25-7f82-4861-9334-3d0e007235bd"
$newGuid = "f19c9bd1-0e77-42d7-b90a-8ffd04108bda"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-09-01 05:04:33"

# Re-create the B2 hyperlink with the refreshed display text (same target).
$overviewLinkAddress = $wsOverview.Hyperlinks.Item(1).Address
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewLinkAddress, "", "", "e2e\$newGuid.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhCnLinkAddress = $wsZhCn.Hyperlinks.Item(1).Address

$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.9ad371bf50a3e12b4f572ecc0c9e0d5803dbceee.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-01 05:04:28"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# Drop both hyperlinks (A2, I2) and re-add only the A2 one.
$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhCnLinkAddress, "", "", "$newGuid.md")

# Column widths for "Latest Target File" / "Latest Handback File" shrink now
# that the columns hold empty values.
$wsZhCn.Columns.Item(9).ColumnWidth = 17.8
$wsZhCn.Columns.Item(10).ColumnWidth = 20.8

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deDeLinkAddress = $wsDeDe.Hyperlinks.Item(1).Address

$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.9ad371bf50a3e12b4f572ecc0c9e0d5803dbceee.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-01 05:04:33"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

# Drop both hyperlinks (A2, I2) and re-add only the A2 one.
$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deDeLinkAddress, "", "", "$newGuid.md")

$wsDeDe.Columns.Item(9).ColumnWidth = 17.8
$wsDeDe.Columns.Item(10).ColumnWidth = 20.8
